$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.505.39"
$ws.Range("E2").Value = "  +5.42%  "
$ws.Range("D3").Value = "2.057.00"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'253.26"
$ws.Range("E5").Value = "  +2.93%  "
$ws.Range("D6").Value = "'0.654"
$ws.Range("E6").Value = "  +3.04%  "
$ws.Range("D7").Value = "'67.79"
$ws.Range("E7").Value = "  +15.84%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.385"
$ws.Range("E9").Value = "  +6.62%  "
$ws.Range("D10").Value = "'59.71"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").Value = "'0.0769"
$ws.Range("E11").Value = "  +4.60%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "'0.930"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "'14.89"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "'22.68"
$ws.Range("E15").Value = "  +24.58%  "
$ws.Range("D16").Value = "2.357.40"
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("D17").Value = "'5.59"
$ws.Range("E17").Value = "  +5.18%  "
$ws.Range("D18").Value = "2.057.64"
$ws.Range("E18").Value = "  +3.72%  "
$ws.Range("D19").Value = "37.397.71"
$ws.Range("E19").Value = "  +5.24%  "
$ws.Range("D20").Value = "'73.70"
$ws.Range("D21").Value = "0.0₃0877"
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("D22").Value = "'5.48"
$ws.Range("E22").Value = "  +4.98%  "
$ws.Range("D23").Value = "'240.17"
$ws.Range("E23").Value = "  +3.30%  "
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  +7.05%  "
$ws.Range("D27").Value = "'9.97"
$ws.Range("E27").Value = "  +9.29%  "
$ws.Range("D28").Value = "'162.25"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "'20.03"
$ws.Range("E29").Value = "  +4.27%  "
$ws.Range("D30").Value = "'0.129"
$ws.Range("E30").Value = "  +35.00%  "
$ws.Range("D31").Value = "'0.123"
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("E32").Value = "  +7.68%  "
$ws.Range("E33").Value = "  +8.79%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.69"
$ws.Range("E34").Value = "  +7.93%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0627"
$ws.Range("E35").Value = "  +5.48%  "
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  +15.34%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("D40").Value = "'3.16"
$ws.Range("E40").Value = "  +39.00%  "
$ws.Range("D41").Value = "'0.104"
$ws.Range("E41").Value = "  +15.39%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.26"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'3.04"
$ws.Range("E43").Value = "  +5.73%  "
$ws.Range("D44").Value = "'17.47"
$ws.Range("E44").Value = "  +7.96%  "
$ws.Range("E45").Value = "  +5.68%  "
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").Value = "'97.33"
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("D48").Value = "'7.95"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").Value = "1.413.78"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("E51").Value = "  +11.30%  "
